# Mantel correlogram table: update the "p" column (4th column) values
# for several distance-class rows, per the recorded diff.
#
# Table layout (1 header row + 14 data rows):
#   col 1 = Distance Class (m), col 2 = N, col 3 = Mantel r, col 4 = p

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($rowIndex, $colIndex, $newValue) {
    $cell = $t.Rows.Item($rowIndex).Cells.Item($colIndex)
    $cellRange = $cell.Range
    # Drop the trailing end-of-cell marker so we only replace the visible text
    $cellRange.End = $cellRange.End - 1
    $cellRange.Text = $newValue
}

# Data row index (2-based, since row 1 is the header) -> new "p" value
Set-CellText 2  4 "0.481"   # 1,250  : 0.461 -> 0.481
Set-CellText 3  4 "0.951"   # 3,750  : 0.909 -> 0.951
Set-CellText 4  4 "0.707"   # 6,250  : 0.647 -> 0.707
Set-CellText 7  4 "1"       # 13,750 : 0.995 -> 1
Set-CellText 8  4 "0.315"   # 16,250 : 0.266 -> 0.315
Set-CellText 9  4 "0.993"   # 18,750 : 1     -> 0.993
Set-CellText 10 4 "0.405"   # 21,250 : 0.315 -> 0.405
Set-CellText 12 4 "0.881"   # 26,250 : 0.89  -> 0.881
Set-CellText 13 4 "0.979"   # 28,750 : 0.919 -> 0.979
Set-CellText 14 4 "1"       # 31,250 : 0.967 -> 1

Write-Output "Updated 10 p-value cells in Mantel correlogram table."
